# Refresh the cryptos price list with the latest scraped values.
# Cells in column D that look like a plain number (e.g. "207.98") are
# forced to stay text (matching the source sheet's inlineStr cells) by
# briefly switching NumberFormat to "@" before the write, then the
# cell style is reset back to Normal so no visible formatting changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.258.83'
$ws.Range("E2").Value = '  -1.42%  '
$ws.Range("D3").Value = '1.574.93'
$ws.Range("E3").Value = '  -0.92%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("E6").Value = '  -1.85%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.26'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -1.16%  '
$ws.Range("E10").Value = '  +0.16%  '
$ws.Range("E11").Value = '  +0.00%  '
$ws.Range("D12").Value = '1.799.09'
$ws.Range("E12").Value = '  -0.92%  '
$ws.Range("D13").Value = '1.580.53'
$ws.Range("E13").Value = '  -0.75%  '
$ws.Range("E14").Value = '  -1.16%  '
$ws.Range("E15").Value = '  -1.30%  '
$ws.Range("D16").Value = '27.270.01'
$ws.Range("E16").Value = '  -1.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.44'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.04%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '214.94'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.35'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.49%  '
$ws.Range("D20").Value = '0.0₃0687'
$ws.Range("E20").Value = '  -0.86%  '
$ws.Range("E21").Value = '  -0.22%  '
$ws.Range("E22").Value = '  -0.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.39'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.75%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.01'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.87'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.87%  '
$ws.Range("E26").Value = '  -4.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.96'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.60%  '
$ws.Range("B28").Value = 'BinanceUSD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.20%  '
$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.104'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.92%  '
$ws.Range("E30").Value = '  -1.68%  '
$ws.Range("E32").Value = '  -1.30%  '
$ws.Range("D33").Value = '1.408.12'
$ws.Range("E33").Value = '  +2.55%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.92'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.47%  '
$ws.Range("E35").Value = '  +1.62%  '
$ws.Range("E36").Value = '  -1.21%  '
$ws.Range("E37").Value = '  -2.72%  '
$ws.Range("E38").Value = '  -1.65%  '
$ws.Range("E39").Value = '  -0.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.518'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.79%  '
$ws.Range("E41").Value = '  -0.16%  '
$ws.Range("E42").Value = '  +2.47%  '
$ws.Range("E43").Value = '  +3.25%  '
$ws.Range("E44").Value = '  +1.99%  '
$ws.Range("B45").Value = 'MXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.18'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.41%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '63.81'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.62%  '
$ws.Range("D47").Value = '1.711.50'
$ws.Range("E47").Value = '  -0.92%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.02'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.15%  '
$ws.Range("E49").Value = '  -1.39%  '
$ws.Range("E50").Value = '  -0.82%  '
$ws.Range("E51").Value = '  -0.04%  '
